$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D27").Value = "루다 선톡을 대비하는법"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/prescaler/"

$ws.Range("D32").Value = "[Airflow] task, dag 우선순위 설정 (priority_weight)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/466"
